$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""58.513.14"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E2").Formula = "=""  -4.12%  """
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D3").Formula = "=""2.537.05"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E3").Formula = "=""  -3.76%  """
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E4").Formula = "=""  -0.04%  """
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D5").Formula = "=""507.89"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E5").Formula = "=""  -4.44%  """
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D6").Formula = "=""144.05"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E6").Formula = "=""  -7.53%  """
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E7").Formula = "=""  +0.00%  """
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D8").Formula = "=""0.566"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E8").Formula = "=""  -4.12%  """
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D9").Formula = "=""2.539.36"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E9").Formula = "=""  -4.04%  """
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E10").Formula = "=""  -7.59%  """
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E11").Formula = "=""  -7.46%  """
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E12").Formula = "=""  -5.46%  """
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E13").Formula = "=""  -0.74%  """
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D14").Formula = "=""2.981.61"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E14").Formula = "=""  -3.87%  """
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D15").Formula = "=""58.501.52"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E15").Formula = "=""  -4.20%  """
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D16").Formula = "=""20.71"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E16").Formula = "=""  -5.90%  """
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E17").Formula = "=""  -6.64%  """
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D18").Formula = "=""2.540.52"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E18").Formula = "=""  -3.73%  """
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E19").Formula = "=""  -5.36%  """
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D20").Formula = "=""334.62"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E20").Formula = "=""  -6.18%  """
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D21").Formula = "=""10.08"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E21").Formula = "=""  -5.58%  """
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D23").Formula = "=""5.95"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E23").Formula = "=""  -4.99%  """
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D24").Formula = "=""60.31"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E24").Formula = "=""  -2.38%  """
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D25").Formula = "=""0.409"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E25").Formula = "=""  -5.33%  """
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D26").Formula = "=""0.999"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E26").Formula = "=""  -0.11%  """
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E27").Formula = "=""  -5.88%  """
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D28").Formula = "=""2.651.91"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E28").Formula = "=""  -3.57%  """
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D29").Formula = "=""0.0₃0785"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E29").Formula = "=""  -9.79%  """
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E30").Formula = "=""  -6.44%  """
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E31").Formula = "=""  +0.00%  """
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D32").Formula = "=""149.57"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E32").Formula = "=""  -1.08%  """
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E33").Formula = "=""  -4.92%  """
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D34").Formula = "=""18.52"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E34").Formula = "=""  -5.32%  """
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E35").Formula = "=""  -5.70%  """
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E36").Formula = "=""  +4.55%  """
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E37").Formula = "=""  -7.01%  """
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E38").Formula = "=""  -7.73%  """
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D39").Formula = "=""36.00"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E39").Formula = "=""  -1.71%  """
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D40").Formula = "=""0.825"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E40").Formula = "=""  -11.07%  """
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E41").Formula = "=""  -7.13%  """
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D42").Formula = "=""283.67"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E42").Formula = "=""  -4.65%  """
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E43").Formula = "=""  -8.06%  """
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D44").Formula = "=""0.0997"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E44").Formula = "=""  -2.50%  """
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E45").Formula = "=""  -0.01%  """
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E46").Formula = "=""  -6.61%  """
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("D47").Formula = "=""0.0534"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E47").Formula = "=""  -5.31%  """
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E48").Formula = "=""  -6.18%  """
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E49").Formula = "=""  -0.44%  """
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E50").Formula = "=""  -5.60%  """
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E51").Formula = "=""  -9.65%  """
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$excel.CutCopyMode = $false
